$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TFC")
$ws2 = $wb.Worksheets.Item("TFC RMSE")

$values = @(43.38389587402344,43.38389587402344,43.38389587402344,43.3838996887207,43.3838996887207,43.38389587402344,43.38389587402344,43.38389587402344,43.38389587402344,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.38389587402344,43.38389587402344,43.38389587402344,43.38389587402344,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.38390731811523,43.38390731811523,43.38390731811523,43.38390731811523,43.38390731811523,43.3838996887207,43.3838996887207,43.3838996887207,43.38390731811523,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.38389587402344,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.38390731811523,43.38390731811523,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.3838996887207,43.38390731811523,43.38390731811523,43.38390731811523,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.38391494750977,43.38391494750977,43.38391494750977,43.3839111328125,43.3839111328125,43.38391494750977,43.38391494750977,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.38391494750977,43.38391494750977,43.38391494750977,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.38391494750977,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.3839111328125,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.3839111328125,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391494750977,43.38391876220703,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38392639160156,43.38393020629883,43.38392639160156,43.38392639160156,43.38392639160156,43.38393020629883,43.38393020629883,43.38392639160156,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38392639160156,43.38392639160156,43.38393020629883,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38393020629883,43.38392639160156,43.38393020629883,43.38392639160156,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393402099609,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393402099609,43.38393402099609,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393402099609,43.38393020629883,43.38393020629883,43.38393402099609,43.38393402099609,43.38393402099609,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38392639160156,43.38391876220703,43.38392639160156,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391494750977,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391494750977,43.38391494750977,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391494750977,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38392639160156,43.38392639160156,43.38392639160156,43.38393020629883,43.38392639160156,43.38392639160156,43.38392639160156,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38391876220703,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38391876220703,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38392639160156,43.38391876220703,43.38391876220703,43.38392639160156,43.38392639160156,43.38392639160156,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393020629883,43.38393402099609,43.38393402099609,43.38393783569336,43.38394165039062,43.38394165039062,43.38394165039062,43.38394927978516)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws1.Cells.Item($i + 2, 2).Value = $values[$i]
}

$ws2.Range("B2").Value = 11.70831566715378
